$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 357.58334
$ws.Range("I2").Value = 230.16667
$ws.Range("J2").Value = 485
$ws.Range("K2").Value = 230.16667
$ws.Range("L2").Value = 485
$ws.Range("M2").Value = -117.16667
$ws.Range("N2").Value = -711

# Row 32
$ws.Range("H32").Value = 1817.5
$ws.Range("I32").Value = 1300.25
$ws.Range("J32").Value = 2024.4
$ws.Range("K32").Value = 1300.25
$ws.Range("L32").Value = 2024.4
$ws.Range("M32").Value = -974.25
$ws.Range("N32").Value = -2676.4

# Row 55
$ws.Range("H55").Value = 210.15384
$ws.Range("I55").Value = 166
$ws.Range("J55").Value = 261.66666
$ws.Range("K55").Value = 166
$ws.Range("L55").Value = 261.66666
$ws.Range("M55").Value = 48
$ws.Range("N55").Value = -689.66666

# Row 103
$ws.Range("H103").Value = 1100
$ws.Range("I103").Value = 773.3333
$ws.Range("K103").Value = 2319.9999
$ws.Range("M103").Value = -1733.9999

# Row 134
$ws.Range("H134").Value = 70780
$ws.Range("J134").Value = 70780
$ws.Range("L134").Value = 70780
$ws.Range("N134").Value = -80920

# Row 137
$ws.Range("H137").Value = 1559.3928
$ws.Range("I137").Value = 1258.3478
$ws.Range("J137").Value = 2944.2
$ws.Range("K137").Value = 3775.0434
$ws.Range("L137").Value = 8832.599999999999
$ws.Range("M137").Value = -1225.0434
$ws.Range("N137").Value = -13932.6

# Row 141
$ws.Range("H141").Value = 1419.4445
$ws.Range("I141").Value = 1419.4445
$ws.Range("K141").Value = 4258.333500000001
$ws.Range("M141").Value = 921.6664999999994

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2370.8289
$ws.Range("I32").Value = 2447.4243
$ws.Range("J32").Value = 1865.3
$ws.Range("K32").Value = 2447.4243
$ws.Range("L32").Value = 1865.3
$ws.Range("M32").Value = -2160.4243
$ws.Range("N32").Value = -2439.3

# Row 61
$ws.Range("H61").Value = 1062.3572
$ws.Range("I61").Value = 881.0833
$ws.Range("J61").Value = 2150
$ws.Range("K61").Value = 881.0833
$ws.Range("L61").Value = 2150
$ws.Range("M61").Value = -669.0833
$ws.Range("N61").Value = -2574

# Row 136
$ws.Range("H136").Value = 1062.3572
$ws.Range("I136").Value = 881.0833
$ws.Range("J136").Value = 2150
$ws.Range("K136").Value = 2643.2499
$ws.Range("L136").Value = 6450
$ws.Range("M136").Value = -93.2498999999998
$ws.Range("N136").Value = -11550

# Row 139
$ws.Range("H139").Value = 47028.75
$ws.Range("J139").Value = 47028.75
$ws.Range("L139").Value = 47028.75
$ws.Range("N139").Value = -57308.75

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 194.5
$ws.Range("I22").Value = 194.5
$ws.Range("K22").Value = 194.5
$ws.Range("M22").Value = -21.5

# Row 94
$ws.Range("H94").Value = 10417379
$ws.Range("I94").Value = 14706420
$ws.Range("J94").Value = 1134.2858
$ws.Range("K94").Value = 14706420
$ws.Range("L94").Value = 1134.2858
$ws.Range("M94").Value = -14705969
$ws.Range("N94").Value = -2036.2858

# Row 107
$ws.Range("H107").Value = 1783.4166
$ws.Range("I107").Value = 1532
$ws.Range("J107").Value = 2202.4443
$ws.Range("K107").Value = 1532
$ws.Range("L107").Value = 2202.4443
$ws.Range("M107").Value = 388
$ws.Range("N107").Value = -6042.4443

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 730.1539
$ws.Range("I58").Value = 730.1539
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 730.1539
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -527.1539
$ws.Range("N58").ClearContents()

# Row 122
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

# Row 132
$ws.Range("H132").Value = 14868.875
$ws.Range("I132").Value = 16564.428
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 49693.284
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -47163.284
$ws.Range("N132").Value = -14060

# Row 134
$ws.Range("H134").Value = 11495306
$ws.Range("I134").Value = 11905763
$ws.Range("K134").Value = 35717289
$ws.Range("M134").Value = -35714754

# Row 136
$ws.Range("H136").Value = 730.1539
$ws.Range("I136").Value = 730.1539
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2190.4617
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 359.5383000000002
$ws.Range("N136").ClearContents()

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 7523.769
$ws.Range("I56").Value = 7523.769
$ws.Range("K56").Value = 7523.769
$ws.Range("M56").Value = -6993.769

# Row 68
$ws.Range("H68").Value = 2294.0264
$ws.Range("J68").Value = 2342.5134
$ws.Range("L68").Value = 7027.540199999999
$ws.Range("N68").Value = -8649.5402

# Row 71
$ws.Range("H71").Value = 2294.0264
$ws.Range("J71").Value = 2342.5134
$ws.Range("L71").Value = 21082.6206
$ws.Range("N71").Value = -29194.6206

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 96
$ws.Range("H96").Value = 7930.3335
$ws.Range("J96").Value = 7930.3335
$ws.Range("L96").Value = 23791.0005
$ws.Range("N96").Value = -27909.0005

# Row 97
$ws.Range("H97").Value = 632.3333
$ws.Range("I97").Value = 632.3333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1896.9999
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1400.9999
$ws.Range("N97").ClearContents()

# Row 107
$ws.Range("H107").Value = 34385.332
$ws.Range("J107").Value = 50801.5
$ws.Range("L107").Value = 152404.5
$ws.Range("N107").Value = -156244.5

# Row 131
$ws.Range("H131").Value = 11237136
$ws.Range("J131").Value = 1224.2073
$ws.Range("L131").Value = 3672.6219
$ws.Range("N131").Value = -13752.6219

# Row 137
$ws.Range("H137").Value = 7898.8965
$ws.Range("I137").Value = 1494.5454
$ws.Range("J137").Value = 11812.667
$ws.Range("K137").Value = 4483.6362
$ws.Range("L137").Value = 35438.001
$ws.Range("M137").Value = 616.3638000000001
$ws.Range("N137").Value = -45638.001

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 916414.75
$ws.Range("I107").Value = 1603150.4
$ws.Range("J107").Value = 767.3333
$ws.Range("K107").Value = 1603150.4
$ws.Range("L107").Value = 767.3333
$ws.Range("M107").Value = -1601230.4
$ws.Range("N107").Value = -4607.3333

# Row 113
$ws.Range("H113").Value = 2368.125
$ws.Range("I113").Value = 1768
$ws.Range("K113").Value = 1768
$ws.Range("M113").Value = 402

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1219.8
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1219.8
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1219.8
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1809.8

# Row 27
$ws.Range("H27").Value = 1219.8
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1219.8
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1219.8
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1433.8

# Row 61
$ws.Range("H61").Value = 3999.3333
$ws.Range("I61").Value = 3499
$ws.Range("K61").Value = 3499
$ws.Range("M61").Value = -3297

# Row 93
$ws.Range("H93").Value = 1125
$ws.Range("I93").Value = 1125
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1125
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 123
$ws.Range("N93").ClearContents()

# Row 113
$ws.Range("H113").Value = 3999.3333
$ws.Range("I113").Value = 3499
$ws.Range("K113").Value = 3499
$ws.Range("M113").Value = -1329

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132
$ws.Range("H132").Value = 33946.71
$ws.Range("I132").Value = 1360
$ws.Range("J132").Value = 145672.58
$ws.Range("K132").Value = 4080
$ws.Range("L132").Value = 437017.74
$ws.Range("M132").Value = -1550
$ws.Range("N132").Value = -442077.74

# Row 136
$ws.Range("H136").Value = 8175.857
$ws.Range("I136").Value = 8710.154
$ws.Range("K136").Value = 26130.462
$ws.Range("M136").Value = -23580.462

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 122
$ws.Range("H122").Value = 17340168
$ws.Range("I122").Value = 28897556
$ws.Range("K122").Value = 86692668
$ws.Range("M122").Value = -86690218

# Row 132
$ws.Range("H132").Value = 1817.7021
$ws.Range("I132").Value = 1628.7949
$ws.Range("J132").Value = 2738.625
$ws.Range("K132").Value = 4886.384700000001
$ws.Range("L132").Value = 8215.875
$ws.Range("M132").Value = -2356.384700000001
$ws.Range("N132").Value = -13275.875

# Row 136
$ws.Range("H136").Value = 631.7143
$ws.Range("I136").Value = 333
$ws.Range("K136").Value = 999
$ws.Range("M136").Value = 1551

Write-Output "All edits applied."